# Edit script: applies the changes described by the diff to before.docx
# Summary of changes:
#  1. Remove the _GoBack bookmark after "BASES TEORICAS" (it gets relocated
#     later, near the REACT paragraph, as a side effect of normal editing -
#     Word moves _GoBack to the last edited spot; we replicate that by
#     re-adding the bookmark at the new location, which removes the old one
#     because bookmark names must be unique).
#  2. Strip <w:lang w:val="en-US"/> from the "4.1 APACHE", "5.2 CSS", the
#     paragraph that becomes "5.4 JAVASCRIPT" and the paragraph that becomes
#     "5.5 REACT".
#  3. Text edits:
#       - "... STYLE SHEETS) " -> "... STYLE SHEETS)" (drop trailing space)
#       - insert a new paragraph "5.3. BOOTSTRAP "
#       - renumber "5.3. JAVASCRIPT" -> "5.4. JAVASCRIPT"
#       - renumber "5.4. REACT" -> "5.5. REACT"
#  4. Re-insert the _GoBack bookmark between the "5" run and the ". REACT"
#     run of the (renumbered) REACT paragraph.

$d = $word.ActiveDocument

function Get-ParaByText($text) {
    $rng = $d.Content
    $found = $rng.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Paragraph containing '$text' not found"
    }
    return $rng.Paragraphs(1)
}

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$pPrNoLang = '<w:pPr><w:widowControl w:val="0"/><w:tabs><w:tab w:val="right" w:leader="dot" w:pos="8261"/></w:tabs><w:spacing w:after="0" w:line="480" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:szCs w:val="22"/></w:rPr></w:pPr>'
$rPrNoLang = '<w:rPr><w:szCs w:val="22"/></w:rPr>'

# --- 1. "4.1. APACHE" paragraph: drop <w:lang> only (text unchanged) ---
$p = Get-ParaByText("APACHE")
$xml = $pkgOpen + '<w:p>' + $pPrNoLang +
    '<w:r>' + $rPrNoLang + '<w:t>4</w:t></w:r>' +
    '<w:r>' + $rPrNoLang + '<w:t>.1.</w:t></w:r>' +
    '<w:r>' + $rPrNoLang + '<w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r>' + $rPrNoLang + '<w:t xml:space="preserve">APACHE </w:t></w:r>' +
    '</w:p>' + $pkgClose
$p.Range.InsertXML($xml)

# --- 2. "5.2. CSS (CASCADE STYLE SHEETS)" paragraph: drop <w:lang>, drop
#        trailing space on the last run ---
$p = Get-ParaByText("CASCAD")
$xml = $pkgOpen + '<w:p>' + $pPrNoLang +
    '<w:r>' + $rPrNoLang + '<w:t>5.</w:t></w:r>' +
    '<w:r>' + $rPrNoLang + '<w:t>2</w:t></w:r>' +
    '<w:r>' + $rPrNoLang + '<w:t>. CSS (CASCAD</w:t></w:r>' +
    '<w:r>' + $rPrNoLang + '<w:t>E</w:t></w:r>' +
    '<w:r>' + $rPrNoLang + '<w:t xml:space="preserve"> STYLE SHEETS)</w:t></w:r>' +
    '</w:p>' + $pkgClose
$p.Range.InsertXML($xml)

# --- 3. Insert new "5.3. BOOTSTRAP " paragraph right after the CSS one ---
$p = Get-ParaByText("CASCAD")
$p.Range.InsertParagraphAfter()
$p = Get-ParaByText("CASCAD")
$newPara = $p.Next(1)
$xml = $pkgOpen + '<w:p>' + $pPrNoLang +
    '<w:r>' + $rPrNoLang + '<w:t>5.3. BOOTSTRAP</w:t></w:r>' +
    '<w:r>' + $rPrNoLang + '<w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p>' + $pkgClose
$newPara.Range.InsertXML($xml)

# --- 4. "5.3. JAVASCRIPT" paragraph -> becomes "5.4. JAVASCRIPT", drop <w:lang> ---
$p = Get-ParaByText("JAVASCRIPT")
$xml = $pkgOpen + '<w:p>' + $pPrNoLang +
    '<w:r>' + $rPrNoLang + '<w:t>5.</w:t></w:r>' +
    '<w:r>' + $rPrNoLang + '<w:t>4</w:t></w:r>' +
    '<w:r>' + $rPrNoLang + '<w:t xml:space="preserve">. JAVASCRIPT </w:t></w:r>' +
    '</w:p>' + $pkgClose
$p.Range.InsertXML($xml)

# --- 5. "5.4. REACT" paragraph -> becomes "5.5. REACT", drop <w:lang> ---
$p = Get-ParaByText("REACT")
$xml = $pkgOpen + '<w:p>' + $pPrNoLang +
    '<w:r>' + $rPrNoLang + '<w:t>5.</w:t></w:r>' +
    '<w:r>' + $rPrNoLang + '<w:t>5</w:t></w:r>' +
    '<w:r>' + $rPrNoLang + '<w:t>. REACT</w:t></w:r>' +
    '</w:p>' + $pkgClose
$p.Range.InsertXML($xml)

# --- 6. Move the _GoBack bookmark: insert between the "5" run and the
#        ". REACT" run of the REACT paragraph. Adding a bookmark with the
#        same name as an existing one moves it (removing the old one). ---
$rng = $d.Content
$found = $rng.Find.Execute("REACT", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "REACT not found for bookmark placement"
}
$bmRange = $d.Range($rng.Start - 1, $rng.Start - 1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Host "done"
